# The sheet currently holds rows for years 2007-2013 (rows 2-8).
# The target keeps only 2010-2013 (shifted up to rows 2-5), i.e. the
# three oldest years (2007, 2008, 2009 -> rows 2:4) are removed and the
# remaining rows shift up, shrinking the used range to A1:F5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("2:4").Delete()
